$p = $ppt.ActivePresentation
$tm = $p.TitleMaster
try {
  $cs = $tm.ColorScheme
  $cs.Colors(1).RGB = 0xABCDEF
  Write-Host "set ok"
} catch {
  Write-Host "ERR:" $_.Exception.Message
}
